$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update Runmode (column C) values: swap Y/N so only 4 modules run
$ws.Range("C2").Value = "N"
$ws.Range("C3").Value = "Y"
$ws.Range("C5").Value = "Y"
$ws.Range("C6").Value = "Y"
$ws.Range("C7").Value = "Y"

# Update selection to C8
$ws.Range("C8").Select()
